$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new quarterly columns (R = 31/03/2024, S = 30/06/2024) by
# copying the formatting of the last existing quarter column (Q) across
# the whole used range, then writing in the new values.
$ws.Range("Q1:Q80").Copy()
$ws.Range("R1:S80").PasteSpecial(-4122)

# Row 1 headers for the two new period columns
$ws.Cells.Item(1, 18).Value = "31/03/2024"
$ws.Cells.Item(1, 19).Value = "30/06/2024"

# Quarterly figures for 31/03/2024 (R) and 30/06/2024 (S)
$ws.Cells.Item(2, 18).Value = 2456846.08
$ws.Cells.Item(2, 19).Value = 2537076.992
$ws.Cells.Item(3, 18).Value = 746212.992
$ws.Cells.Item(3, 19).Value = 728403.968
$ws.Cells.Item(4, 18).Value = 188668
$ws.Cells.Item(4, 19).Value = 251338
$ws.Cells.Item(5, 18).Value = 38035
$ws.Cells.Item(5, 19).Value = 51212
$ws.Cells.Item(6, 18).Value = 426723.008
$ws.Cells.Item(6, 19).Value = 305188.992
$ws.Cells.Item(7, 18).Value = 6052
$ws.Cells.Item(7, 19).Value = 5198
$ws.Cells.Item(8, 18).Value = 0
$ws.Cells.Item(8, 19).Value = 0
$ws.Cells.Item(9, 18).Value = 56125
$ws.Cells.Item(9, 19).Value = 76806
$ws.Cells.Item(10, 18).Value = 18323
$ws.Cells.Item(10, 19).Value = 14505
$ws.Cells.Item(11, 18).Value = 12287
$ws.Cells.Item(11, 19).Value = 24156
$ws.Cells.Item(12, 18).Value = 229846
$ws.Cells.Item(12, 19).Value = 217374
$ws.Cells.Item(13, 18).Value = 0
$ws.Cells.Item(13, 19).Value = 0
$ws.Cells.Item(14, 18).Value = 0
$ws.Cells.Item(14, 19).Value = 0
$ws.Cells.Item(15, 18).Value = 9724
$ws.Cells.Item(15, 19).Value = 10182
$ws.Cells.Item(16, 18).Value = 0
$ws.Cells.Item(16, 19).Value = 0
$ws.Cells.Item(17, 18).Value = 0
$ws.Cells.Item(17, 19).Value = 0
$ws.Cells.Item(18, 18).Value = 0
$ws.Cells.Item(18, 19).Value = 0
$ws.Cells.Item(19, 18).Value = 165588.992
$ws.Cells.Item(19, 19).Value = 152128
$ws.Cells.Item(20, 18).Value = 0
$ws.Cells.Item(20, 19).Value = 166
$ws.Cells.Item(21, 18).Value = 0
$ws.Cells.Item(21, 19).Value = 0
$ws.Cells.Item(22, 18).Value = 0
$ws.Cells.Item(22, 19).Value = 0
$ws.Cells.Item(23, 18).Value = 1466139.008
$ws.Cells.Item(23, 19).Value = 1575522.944
$ws.Cells.Item(24, 18).Value = 14648
$ws.Cells.Item(24, 19).Value = 15776
$ws.Cells.Item(25, 18).Value = 0
$ws.Cells.Item(25, 19).Value = 0
$ws.Cells.Item(26, 18).Value = 2456846.08
$ws.Cells.Item(26, 19).Value = 2537076.992
$ws.Cells.Item(27, 18).Value = 622953.9840000001
$ws.Cells.Item(27, 19).Value = 625985.024
$ws.Cells.Item(28, 18).Value = 112009
$ws.Cells.Item(28, 19).Value = 109922
$ws.Cells.Item(29, 18).Value = 138852.992
$ws.Cells.Item(29, 19).Value = 80681
$ws.Cells.Item(30, 18).Value = 28867
$ws.Cells.Item(30, 19).Value = 24608
$ws.Cells.Item(31, 18).Value = 298872.992
$ws.Cells.Item(31, 19).Value = 353313.984
$ws.Cells.Item(32, 18).Value = 0
$ws.Cells.Item(32, 19).Value = 0
$ws.Cells.Item(33, 18).Value = 0
$ws.Cells.Item(33, 19).Value = 0
$ws.Cells.Item(34, 18).Value = 44352
$ws.Cells.Item(34, 19).Value = 57460
$ws.Cells.Item(35, 18).Value = 0
$ws.Cells.Item(35, 19).Value = 0
$ws.Cells.Item(36, 18).Value = 0
$ws.Cells.Item(36, 19).Value = 0
$ws.Cells.Item(37, 18).Value = 1001251.008
$ws.Cells.Item(37, 19).Value = 998876.032
$ws.Cells.Item(38, 18).Value = 890603.008
$ws.Cells.Item(38, 19).Value = 913593.024
$ws.Cells.Item(39, 18).Value = 0
$ws.Cells.Item(39, 19).Value = 0
$ws.Cells.Item(40, 18).Value = 80710
$ws.Cells.Item(40, 19).Value = 73475
$ws.Cells.Item(41, 18).Value = 1078
$ws.Cells.Item(41, 19).Value = 1199
$ws.Cells.Item(42, 18).Value = 0
$ws.Cells.Item(42, 19).Value = 0
$ws.Cells.Item(43, 18).Value = 28860
$ws.Cells.Item(43, 19).Value = 10609
$ws.Cells.Item(44, 18).Value = 0
$ws.Cells.Item(44, 19).Value = 0
$ws.Cells.Item(45, 18).Value = 0
$ws.Cells.Item(45, 19).Value = 0
$ws.Cells.Item(46, 18).Value = 0
$ws.Cells.Item(46, 19).Value = 0
$ws.Cells.Item(47, 18).Value = 832641.024
$ws.Cells.Item(47, 19).Value = 912216
$ws.Cells.Item(48, 18).Value = 803662.976
$ws.Cells.Item(48, 19).Value = 803662.976
$ws.Cells.Item(49, 18).Value = 87694
$ws.Cells.Item(49, 19).Value = 87989
$ws.Cells.Item(50, 18).Value = 0
$ws.Cells.Item(50, 19).Value = 0
$ws.Cells.Item(51, 18).Value = 1443
$ws.Cells.Item(51, 19).Value = 1047
$ws.Cells.Item(52, 18).Value = -48359
$ws.Cells.Item(52, 19).Value = -53052
$ws.Cells.Item(53, 18).Value = -217
$ws.Cells.Item(53, 19).Value = -3198
$ws.Cells.Item(54, 18).Value = 0
$ws.Cells.Item(54, 19).Value = 0
$ws.Cells.Item(55, 18).Value = -11583
$ws.Cells.Item(55, 19).Value = 75767
$ws.Cells.Item(56, 18).Value = 0
$ws.Cells.Item(56, 19).Value = 0
$ws.Cells.Item(59, 18).Value = 490028
$ws.Cells.Item(59, 19).Value = 383608.992
$ws.Cells.Item(60, 18).Value = -354999.008
$ws.Cells.Item(60, 19).Value = -283115.008
$ws.Cells.Item(61, 18).Value = 135028.992
$ws.Cells.Item(61, 19).Value = 100494
$ws.Cells.Item(62, 18).Value = 0
$ws.Cells.Item(62, 19).Value = 0
$ws.Cells.Item(63, 18).Value = -43144
$ws.Cells.Item(63, 19).Value = -44053
$ws.Cells.Item(64, 18).Value = 0
$ws.Cells.Item(64, 19).Value = 0
$ws.Cells.Item(65, 18).Value = 394
$ws.Cells.Item(65, 19).Value = 26450
$ws.Cells.Item(66, 18).Value = -4803
$ws.Cells.Item(66, 19).Value = 0
$ws.Cells.Item(67, 18).Value = 0
$ws.Cells.Item(67, 19).Value = 0
$ws.Cells.Item(68, 18).Value = -47132
$ws.Cells.Item(68, 19).Value = -76768
$ws.Cells.Item(69, 18).Value = 7945
$ws.Cells.Item(69, 19).Value = 8337
$ws.Cells.Item(70, 18).Value = -55077
$ws.Cells.Item(70, 19).Value = -85105
$ws.Cells.Item(74, 18).Value = 40344
$ws.Cells.Item(74, 19).Value = 6123
$ws.Cells.Item(75, 18).Value = -11504
$ws.Cells.Item(75, 19).Value = 355
$ws.Cells.Item(76, 18).Value = -7703
$ws.Cells.Item(76, 19).Value = -11171
$ws.Cells.Item(79, 18).Value = 0
$ws.Cells.Item(79, 19).Value = 0
$ws.Cells.Item(80, 18).Value = 21137
$ws.Cells.Item(80, 19).Value = -4693
